$d = $word.ActiveDocument

# Replace "easily retrieved" with "easy to read"
$d.Content.Find.Execute("easily retrieved", $true, $false, $false, $false, $false,
                         $true, 1, $false, "easy to read", 2)

# Replace "such when" with "such information as when"
$d.Content.Find.Execute("such when", $true, $false, $false, $false, $false,
                         $true, 1, $false, "such information as when", 2)
